# MENAMBAHKAN +62 PADA MENU SISWA
# Applies the recorded edits to data_beasiswa.xlsx:
#  - Siswa: fix Ahmad Maulidun's phone-area note (trailing space on name),
#           move his address to Subang, and register a new student
#           (Ghisya Adi) who is also from Subang.
#  - Beasiswa: Djarum quota freed up by 1 slot and is available again.
#  - Pemberian: remove Deandra's Djarum scholarship grant record (revoked).
#  - Histori_Pencabutan: log Deandra's Djarum scholarship revocation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Siswa"
# ---------------------------------------------------------------
$siswa = $wb.Worksheets.Item("Siswa")

# Row 2: Ahmad Maulidun -> trailing space in name, address Compreng -> Subang
$siswa.Range("B2").Value = "Ahmad Maulidun "
$siswa.Range("D2").Value = "Subang"

# New row 5: Ghisya Adi
$siswa.Range("A5").NumberFormat = "@"
$siswa.Range("A5").Value = "00900980"
$siswa.Range("A5").Style = "Normal"

$siswa.Range("B5").Value = "Ghisya Adi"

$siswa.Range("C5").NumberFormat = "@"
$siswa.Range("C5").Value = "11221111111"
$siswa.Range("C5").Style = "Normal"

$siswa.Range("D5").Value = "Subang"

# ---------------------------------------------------------------
# Sheet "Beasiswa"
# ---------------------------------------------------------------
$beasiswa = $wb.Worksheets.Item("Beasiswa")

# Row 3 (Djarum): quota 0 -> 1, status Habis -> Tersedia
$beasiswa.Range("F3").Value = 1
$beasiswa.Range("G3").Value = "Tersedia"

# ---------------------------------------------------------------
# Sheet "Pemberian" - remove Deandra / B02001 grant (row 3)
# ---------------------------------------------------------------
$pemberian = $wb.Worksheets.Item("Pemberian")
$pemberian.Rows(3).Delete()

# ---------------------------------------------------------------
# Sheet "Histori_Pencabutan" - add revocation record (row 3)
# ---------------------------------------------------------------
$histori = $wb.Worksheets.Item("Histori_Pencabutan")

$histori.Range("A3").NumberFormat = "@"
$histori.Range("A3").Value = "0087654321"
$histori.Range("A3").Style = "Normal"

$histori.Range("B3").Value = "B02001"

$histori.Range("C3").NumberFormat = "@"
$histori.Range("C3").Value = "2025-11-28"
$histori.Range("C3").Style = "Normal"
